# Fix NaN-exclusion bug: p-values/counts recomputed so a single NaN
# behavior event is excluded rather than NaNing out the whole rows p-value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("H2").Value = 0.35714285714285715
$ws.Range("I2").Value = 15.0
$ws.Range("J2").Value = 0.0
$ws.Range("K2").Value = 0.0

# Row 3
$ws.Range("H3").Value = 0.20833333333333334
$ws.Range("I3").Value = 5.0

# Row 6
$ws.Range("F6").Value = 0.11904761904761904
$ws.Range("G6").Value = 5.0
$ws.Range("L6").Value = 0.09523809523809523
$ws.Range("M6").Value = 4.0
$ws.Range("N6").Value = 0.14285714285714285
$ws.Range("O6").Value = 6.0

# Row 8
$ws.Range("J8").Value = 0.057692307692307696
$ws.Range("K8").Value = 3.0

# Row 11
$ws.Range("J11").Value = 0.07894736842105263
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 0.18421052631578946
$ws.Range("M11").Value = 7.0

# Row 12
$ws.Range("D12").Value = 0.15384615384615385
$ws.Range("E12").Value = 6.0
$ws.Range("N12").Value = 0.41025641025641024
$ws.Range("O12").Value = 16.0

# Row 13
$ws.Range("L13").Value = 0.17857142857142858
$ws.Range("M13").Value = 5.0

# Row 14
$ws.Range("N14").Value = 0.375
$ws.Range("O14").Value = 3.0

# Row 17
$ws.Range("F17").Value = 0.3125
$ws.Range("G17").Value = 5.0
$ws.Range("H17").Value = 0.4375
$ws.Range("I17").Value = 7.0

# Row 19
$ws.Range("H19").Value = 0.39215686274509803
$ws.Range("I19").Value = 20.0
$ws.Range("N19").Value = 0.23529411764705882
$ws.Range("O19").Value = 12.0

# Row 20
$ws.Range("N20").Value = 0.11538461538461539
$ws.Range("O20").Value = 6.0

# Row 21
$ws.Range("D21").Value = 0.21428571428571427
$ws.Range("E21").Value = 9.0
$ws.Range("F21").Value = 0.30952380952380953
$ws.Range("G21").Value = 13.0
$ws.Range("L21").Value = 0.0
$ws.Range("M21").Value = 0.0

# Row 22
$ws.Range("H22").Value = 0.23255813953488372
$ws.Range("I22").Value = 10.0
$ws.Range("L22").Value = 0.06976744186046512
$ws.Range("M22").Value = 3.0

# Row 24
$ws.Range("J24").Value = 0.0
$ws.Range("K24").Value = 0.0
$ws.Range("N24").Value = 0.1282051282051282
$ws.Range("O24").Value = 5.0

# Row 26
$ws.Range("H26").Value = 0.41379310344827586
$ws.Range("I26").Value = 12.0

# Row 28
$ws.Range("D28").Value = 0.11538461538461539
$ws.Range("E28").Value = 6.0
$ws.Range("J28").Value = 0.19230769230769232
$ws.Range("K28").Value = 10.0
$ws.Range("L28").Value = 0.23076923076923078
$ws.Range("M28").Value = 12.0

# Row 29
$ws.Range("F29").Value = 0.06666666666666667
$ws.Range("G29").Value = 3.0
$ws.Range("H29").Value = 0.17777777777777778
$ws.Range("I29").Value = 8.0
$ws.Range("J29").Value = 0.15555555555555556
$ws.Range("K29").Value = 7.0
$ws.Range("N29").Value = 0.28888888888888886
$ws.Range("O29").Value = 13.0

# Row 30
$ws.Range("F30").Value = 0.14285714285714285
$ws.Range("G30").Value = 6.0
$ws.Range("L30").Value = 0.16666666666666666
$ws.Range("M30").Value = 7.0

# Row 31
$ws.Range("F31").Value = 0.14705882352941177
$ws.Range("G31").Value = 5.0
$ws.Range("H31").Value = 0.29411764705882354
$ws.Range("I31").Value = 10.0

# Row 32
$ws.Range("H32").Value = 0.0967741935483871
$ws.Range("I32").Value = 6.0
$ws.Range("L32").Value = 0.03225806451612903
$ws.Range("M32").Value = 2.0
$ws.Range("N32").Value = 0.06451612903225806
$ws.Range("O32").Value = 4.0

# Row 34
$ws.Range("H34").Value = 0.20512820512820512
$ws.Range("I34").Value = 8.0
$ws.Range("J34").Value = 0.02564102564102564
$ws.Range("K34").Value = 1.0
$ws.Range("L34").Value = 0.05128205128205128
$ws.Range("M34").Value = 2.0

# Row 36
$ws.Range("D36").Value = 0.25
$ws.Range("E36").Value = 5.0
$ws.Range("F36").Value = 0.45
$ws.Range("G36").Value = 9.0
$ws.Range("J36").Value = 0.05
$ws.Range("K36").Value = 1.0
$ws.Range("N36").Value = 0.1
$ws.Range("O36").Value = 2.0

# Row 40
$ws.Range("F40").Value = 0.09375
$ws.Range("G40").Value = 3.0
$ws.Range("L40").Value = 0.125
$ws.Range("M40").Value = 4.0

# Row 41
$ws.Range("N41").Value = 0.04878048780487805
$ws.Range("O41").Value = 2.0

# Row 43
$ws.Range("D43").Value = 0.03125
$ws.Range("E43").Value = 3.0
$ws.Range("F43").Value = 0.07291666666666667
$ws.Range("G43").Value = 7.0
$ws.Range("J43").Value = 0.052083333333333336
$ws.Range("K43").Value = 5.0
$ws.Range("N43").Value = 0.16666666666666666
$ws.Range("O43").Value = 16.0

# Row 44
$ws.Range("D44").Value = 0.08
$ws.Range("E44").Value = 2.0
$ws.Range("H44").Value = 0.32
$ws.Range("I44").Value = 8.0

# Row 45
$ws.Range("D45").Value = 0.35714285714285715
$ws.Range("E45").Value = 10.0
$ws.Range("F45").Value = 0.42857142857142855
$ws.Range("G45").Value = 12.0
$ws.Range("H45").Value = 0.5
$ws.Range("I45").Value = 14.0
$ws.Range("J45").Value = 0.07142857142857142
$ws.Range("K45").Value = 2.0

# Row 47
$ws.Range("F47").Value = 0.12121212121212122
$ws.Range("G47").Value = 4.0
$ws.Range("N47").Value = 0.06060606060606061
$ws.Range("O47").Value = 2.0

# Row 48
$ws.Range("H48").Value = 0.10256410256410256
$ws.Range("I48").Value = 4.0
$ws.Range("J48").Value = 0.02564102564102564
$ws.Range("K48").Value = 1.0
$ws.Range("N48").Value = 0.1282051282051282
$ws.Range("O48").Value = 5.0

# Row 49
$ws.Range("J49").Value = 0.05357142857142857
$ws.Range("K49").Value = 3.0

# Row 50
$ws.Range("N50").Value = 0.10638297872340426
$ws.Range("O50").Value = 5.0

# Row 51
$ws.Range("F51").Value = 0.08888888888888889
$ws.Range("G51").Value = 4.0
$ws.Range("L51").Value = 0.08888888888888889
$ws.Range("M51").Value = 4.0
$ws.Range("N51").Value = 0.15555555555555556
$ws.Range("O51").Value = 7.0

# Row 52
$ws.Range("H52").Value = 0.075
$ws.Range("I52").Value = 3.0
$ws.Range("N52").Value = 0.075
$ws.Range("O52").Value = 3.0

# Row 54
$ws.Range("D54").Value = 0.02631578947368421
$ws.Range("E54").Value = 1.0

# Row 55
$ws.Range("F55").Value = 0.034482758620689655
$ws.Range("G55").Value = 1.0
$ws.Range("H55").Value = 0.13793103448275862
$ws.Range("I55").Value = 4.0
$ws.Range("L55").Value = 0.10344827586206896
$ws.Range("M55").Value = 3.0
$ws.Range("N55").Value = 0.13793103448275862
$ws.Range("O55").Value = 4.0
